# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rebuilds the worker/period detail table (rows 16-48) so that the new
# employee (NELSON LUIS PEREA ANDRADE, CC 12917341) is interleaved with the
# existing employee (EDGAR LUIS ALMAGRO MENDOZA, CC 1047472256) period by
# period, in ascending period order, with NELSON's short first period
# (2208, partial value 12000) at the top and both employees' final period
# (2312) carrying the prorated value 25333 at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edgarDoc = "1047472256"
$edgarName = "EDGAR LUIS ALMAGRO MENDOZA"
$nelsonDoc = "12917341"
$nelsonName = "NELSON LUIS PEREA ANDRADE"

$rows = @(
    @{Row=16; Doc=$nelsonDoc; Name=$nelsonName; Period="2208"; Valor=12000},
    @{Row=17; Doc=$edgarDoc;  Name=$edgarName;  Period="2209"; Valor=40000},
    @{Row=18; Doc=$nelsonDoc; Name=$nelsonName; Period="2209"; Valor=40000},
    @{Row=19; Doc=$edgarDoc;  Name=$edgarName;  Period="2210"; Valor=40000},
    @{Row=20; Doc=$nelsonDoc; Name=$nelsonName; Period="2210"; Valor=40000},
    @{Row=21; Doc=$edgarDoc;  Name=$edgarName;  Period="2211"; Valor=40000},
    @{Row=22; Doc=$nelsonDoc; Name=$nelsonName; Period="2211"; Valor=40000},
    @{Row=23; Doc=$edgarDoc;  Name=$edgarName;  Period="2212"; Valor=40000},
    @{Row=24; Doc=$nelsonDoc; Name=$nelsonName; Period="2212"; Valor=40000},
    @{Row=25; Doc=$edgarDoc;  Name=$edgarName;  Period="2301"; Valor=40000},
    @{Row=26; Doc=$nelsonDoc; Name=$nelsonName; Period="2301"; Valor=40000},
    @{Row=27; Doc=$edgarDoc;  Name=$edgarName;  Period="2302"; Valor=40000},
    @{Row=28; Doc=$nelsonDoc; Name=$nelsonName; Period="2302"; Valor=40000},
    @{Row=29; Doc=$edgarDoc;  Name=$edgarName;  Period="2303"; Valor=40000},
    @{Row=30; Doc=$nelsonDoc; Name=$nelsonName; Period="2303"; Valor=40000},
    @{Row=31; Doc=$edgarDoc;  Name=$edgarName;  Period="2304"; Valor=40000},
    @{Row=32; Doc=$nelsonDoc; Name=$nelsonName; Period="2304"; Valor=40000},
    @{Row=33; Doc=$edgarDoc;  Name=$edgarName;  Period="2305"; Valor=40000},
    @{Row=34; Doc=$nelsonDoc; Name=$nelsonName; Period="2305"; Valor=40000},
    @{Row=35; Doc=$edgarDoc;  Name=$edgarName;  Period="2306"; Valor=40000},
    @{Row=36; Doc=$nelsonDoc; Name=$nelsonName; Period="2306"; Valor=40000},
    @{Row=37; Doc=$edgarDoc;  Name=$edgarName;  Period="2307"; Valor=40000},
    @{Row=38; Doc=$nelsonDoc; Name=$nelsonName; Period="2307"; Valor=40000},
    @{Row=39; Doc=$edgarDoc;  Name=$edgarName;  Period="2308"; Valor=40000},
    @{Row=40; Doc=$nelsonDoc; Name=$nelsonName; Period="2308"; Valor=40000},
    @{Row=41; Doc=$edgarDoc;  Name=$edgarName;  Period="2309"; Valor=40000},
    @{Row=42; Doc=$nelsonDoc; Name=$nelsonName; Period="2309"; Valor=40000},
    @{Row=43; Doc=$edgarDoc;  Name=$edgarName;  Period="2310"; Valor=40000},
    @{Row=44; Doc=$nelsonDoc; Name=$nelsonName; Period="2310"; Valor=40000},
    @{Row=45; Doc=$edgarDoc;  Name=$edgarName;  Period="2311"; Valor=40000},
    @{Row=46; Doc=$nelsonDoc; Name=$nelsonName; Period="2311"; Valor=40000},
    @{Row=47; Doc=$edgarDoc;  Name=$edgarName;  Period="2312"; Valor=25333},
    @{Row=48; Doc=$nelsonDoc; Name=$nelsonName; Period="2312"; Valor=25333}
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value2 = "CC"
    $ws.Cells.Item($r, 3).Value2 = $item.Doc
    $ws.Cells.Item($r, 4).Value2 = $item.Name
    $ws.Cells.Item($r, 5).Value2 = $item.Period
    $ws.Cells.Item($r, 6).Value2 = $item.Valor
    $ws.Cells.Item($r, 7).Value2 = 1000000
}

"ok"
